# "add to cart feature"
# Fill in the new "add_item_cart" API documentation row (row 17), plus the
# single-letter placeholder "p" that was dropped into B16, matching the
# author's added shared-strings entries (#51-#54) and the resulting sheet
# selection / row height changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New content (row 16 / 17, columns B-D)
$ws.Range("B16").Value = "p"
$ws.Range("B17").Value = "POST`n/add_item_cart"
$ws.Range("C17").Value = "request{`nitem_id,`namount,`n}response{`nall cart items with the cart data "
$ws.Range("D17").Value = "get all cart items for thet client and add it to session `nif the cart has this item so update if not do and update session and return`nsame thing with delete if it exsist in the session so update if not so delete after order is placed delet cart items and  create new cart id`nq: what happpends if somone on another site logs in and changes things around so we need to have a session saver a: what we need to do is save a session to a user on the data base on delete it on log out) "

# Same wrap-text style ("s=4" in the original workbook) as the rest of
# column B's documentation cells.
$ws.Range("B16").WrapText = $true
$ws.Range("B17:D17").WrapText = $true

# Row 17 grew tall to fit the new request/response notes.
$ws.Rows("17").RowHeight = 210

# Author ended up with D17 selected, scrolled so row 9 is at the top.
[void]$ws.Range("D17").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 9
$win.ScrollColumn = 2
